$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.5610447226247582
$ws.Range("C4").Value = 0.534
$ws.Range("D4").Value = 0.5985859442797988
$ws.Range("E4").Value = 0.5860000000000001
$ws.Range("F4").Value = 0.6579176836896354
$ws.Range("G4").Value = 0.89
$ws.Range("H4").Value = 0.5314183833503119
$ws.Range("I4").Value = 0.5385000000000001
$ws.Range("J4").Value = 0.7290315095519146
$ws.Range("K4").Value = 0.849
$ws.Range("L4").Value = 0.6483729490791019
$ws.Range("M4").Value = 0.6829999999999999

$ws.Range("B5").Value = 0.5495558523323343
$ws.Range("C5").Value = 0.448
$ws.Range("D5").Value = 0.7314451844686448
$ws.Range("E5").Value = 0.6389999999999999
$ws.Range("F5").Value = 0.535181020306055
$ws.Range("G5").Value = 0.5309999999999999
$ws.Range("H5").Value = 0.550865745797552
$ws.Range("I5").Value = 0.5595
$ws.Range("J5").Value = 0.6899196395594889
$ws.Range("K5").Value = 0.678
$ws.Range("L5").Value = 0.7249144224478739
$ws.Range("M5").Value = 0.7175

$ws.Range("B6").Value = 0.6416973071959204
$ws.Range("C6").Value = 0.6419999999999999
$ws.Range("D6").Value = 0.6508146326546903
$ws.Range("E6").Value = 0.645
$ws.Range("F6").Value = 0.5501052135271551
$ws.Range("G6").Value = 0.53
$ws.Range("H6").Value = 0.5920101700508398
$ws.Range("I6").Value = 0.5814999999999999
$ws.Range("J6").Value = 0.7419305070373569
$ws.Range("K6").Value = 0.749
$ws.Range("L6").Value = 0.7603074785133486
$ws.Range("M6").Value = 0.744
